$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Coin "Price" values look numeric (e.g. "305.71") and Excel will
    # silently coerce a plain .Value assignment into a floating point
    # number, losing the original text representation (and precision
    # for thousands-grouped values like "42.627.24"). Forcing the cell
    # to Text format for the assignment keeps it a string, then clearing
    # formats afterwards restores the original (unstyled) cell format.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.627.24"
$ws.Range("E2").Value = "  -1.05%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.287.85"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "305.71"
$ws.Range("E5").Value = "  +1.76%  "

# Row 6 - Solana
Set-TextValue "D6" "95.79"
$ws.Range("E6").Value = "  -2.61%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.64%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.07%  "

# Row 10 - Avalanche
Set-TextValue "D10" "35.07"
$ws.Range("E10").Value = "  -3.10%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.37%  "

# Row 12 - Chainlink
Set-TextValue "D12" "18.36"
$ws.Range("E12").Value = "  +3.56%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.01%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.71"
$ws.Range("E14").Value = "  -2.10%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.644.98"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.293.69"
$ws.Range("E16").Value = "  -1.10%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.777"
$ws.Range("E17").Value = "  -1.33%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.549.57"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("E19").Value = "  +0.43%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -2.21%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.03"
$ws.Range("E21").Value = "  -1.44%  "

# Row 22 - Litecoin
Set-TextValue "D22" "66.90"
$ws.Range("E22").Value = "  -3.09%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "234.94"
$ws.Range("E23").Value = "  -0.99%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -0.44%  "

# Row 25 - was PancakeSwap, becomes Dai
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - was Dai, becomes PancakeSwap
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "2.45"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27 - LEO
Set-TextValue "D27" "4.02"
$ws.Range("E27").Value = "  +0.02%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "24.96"
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - was Toncoin, becomes Monero
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D29" "165.72"
$ws.Range("E29").Value = "  +0.56%  "

# Row 30 - was Monero, becomes Toncoin
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "2.05"
$ws.Range("E30").Value = "  +0.40%  "

# Row 31 - Cosmos
Set-TextValue "D31" "9.01"
$ws.Range("E31").Value = "  -1.06%  "

# Row 32 - InjectiveProtocol
Set-TextValue "D32" "32.77"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.07%  "

# Row 34 - RenderToken
Set-TextValue "D34" "4.71"
$ws.Range("E34").Value = "  -1.09%  "

# Row 35 - Filecoin
Set-TextValue "D35" "4.95"
$ws.Range("E35").Value = "  -2.48%  "

# Row 36 - Celestia
$ws.Range("E36").Value = "  -2.28%  "

# Row 37 - WEMIXToken
Set-TextValue "D37" "2.39"
$ws.Range("E37").Value = "  -0.80%  "

# Row 38 - Hedera
Set-TextValue "D38" "0.0690"
$ws.Range("E38").Value = "  -0.88%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.99%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -2.15%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -1.69%  "

# Row 42 - LidoDAOToken
Set-TextValue "D42" "2.68"
$ws.Range("E42").Value = "  -3.50%  "

# Row 43 - Maker
Set-TextValue "D43" "1.994.87"
$ws.Range("E43").Value = "  -0.83%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -3.01%  "

# Row 45 - EnergySwap
$ws.Range("E45").Value = "  +3.82%  "

# Row 46 - FraxShare
Set-TextValue "D46" "10.02"
$ws.Range("E46").Value = "  -2.93%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  -9.58%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  -1.83%  "

# Row 49 - HuobiToken
Set-TextValue "D49" "2.87"
$ws.Range("E49").Value = "  +8.10%  "

# Row 50 - MultiversX
Set-TextValue "D50" "53.64"
$ws.Range("E50").Value = "  -0.79%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.511.47"
$ws.Range("E51").Value = "  -0.48%  "
